# "Generate Report for Handback" — refresh the handoff/handback timestamps
# for the 36014684-8719-4fb6-99ad-182db142a162.md file row across the
# Overview / zh-cn / de-de report sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G), row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-21 04:54:04"

# --- zh-cn sheet: row 3 is the 36014684-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-21 04:53:58"   # Correspond Handoff Datetime
$wsZhCn.Range("K3").Value = "2016-08-21 04:54:26"   # Correspond Handback DateTime

# --- de-de sheet: row 3 is the 36014684-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-21 04:54:04"   # Correspond Handoff Datetime
$wsDeDe.Range("K3").Value = "2016-08-21 04:54:32"   # Correspond Handback DateTime
